# Apply the updated cryptos price/volume snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.903.35"
$ws.Range("E2").Value = "  +0.83%  "

$ws.Range("D3").Value = "2.361.22"
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "0.694"
$ws.Range("E5").Value = "  +5.10%  "

$ws.Range("D6").Value = "242.83"
$ws.Range("E6").Value = "  +3.69%  "

$ws.Range("D7").Value = "77.09"
$ws.Range("E7").Value = "  +4.81%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "0.637"
$ws.Range("E9").Value = "  +23.35%  "

$ws.Range("D10").Value = "0.103"
$ws.Range("E10").Value = "  +4.97%  "

$ws.Range("D11").Value = "57.46"
$ws.Range("E11").Value = "  +0.74%  "

$ws.Range("D12").Value = "34.02"
$ws.Range("E12").Value = "  +23.91%  "

$ws.Range("E13").Value = "  +18.01%  "

$ws.Range("E14").Value = "  +2.04%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "17.12"
$ws.Range("E15").Value = "  +4.36%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.711.38"
$ws.Range("E16").Value = "  -0.44%  "

$ws.Range("D17").Value = "0.933"
$ws.Range("E17").Value = "  +6.55%  "

$ws.Range("D18").Value = "2.360.52"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("D19").Value = "43.795.84"
$ws.Range("E19").Value = "  +0.55%  "

$ws.Range("D20").Value = "0.0000104"
$ws.Range("E20").Value = "  +3.04%  "

$ws.Range("D21").Value = "6.72"
$ws.Range("E21").Value = "  +4.93%  "

$ws.Range("D22").Value = "78.01"
$ws.Range("E22").Value = "  +3.04%  "

$ws.Range("D23").Value = "257.64"
$ws.Range("E23").Value = "  +2.35%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "2.55"
$ws.Range("E25").Value = "  +2.78%  "

$ws.Range("D26").Value = "11.11"
$ws.Range("E26").Value = "  +9.09%  "

$ws.Range("D27").Value = "3.63"
$ws.Range("E27").Value = "  -3.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +17.59%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "23.09"
$ws.Range("E29").Value = "  +2.48%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").Value = "  -0.33%  "

$ws.Range("D31").Value = "175.44"
$ws.Range("E31").Value = "  +2.05%  "

$ws.Range("E32").Value = "  -3.62%  "

$ws.Range("E33").Value = "  +5.26%  "

$ws.Range("D34").Value = "5.38"
$ws.Range("E34").Value = "  +5.32%  "

$ws.Range("D35").Value = "0.0764"
$ws.Range("E35").Value = "  +9.51%  "

$ws.Range("E36").Value = "  +5.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.66%  "

$ws.Range("D38").Value = "2.45"
$ws.Range("E38").Value = "  +0.53%  "

$ws.Range("D39").Value = "6.49"
$ws.Range("E39").Value = "  -2.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0280"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.62%  "

$ws.Range("D41").Value = "19.54"
$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").Value = "0.205"
$ws.Range("E42").Value = "  +16.50%  "

$ws.Range("D43").Value = "9.09"
$ws.Range("E43").Value = "  +2.63%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.105"
$ws.Range("E44").Value = "  +9.60%  "

$ws.Range("B45").Value = "BinanceUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("D46").Value = "2.55"
$ws.Range("E46").Value = "  +12.66%  "

$ws.Range("E47").Value = "  +4.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.48%  "

$ws.Range("D49").Value = "102.47"
$ws.Range("E49").Value = "  +2.33%  "

$ws.Range("D50").Value = "57.24"
$ws.Range("E50").Value = "  +12.58%  "

$ws.Range("D51").Value = "4.49"
$ws.Range("E51").Value = "  -0.07%  "
